# Scen_NCAP_NUC.xlsx - update SMR / wind-on / wind-off upper-bound values
# and restore the sheet view (zoom/selection) left behind by the author.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- ELE_NEW_NUC_SMR / "Ograniczenia gorne dla SMR" block (rows 15-20, 22) ---
$ws.Range("E15").Value = 10
$ws.Range("E16").Value = 10
$ws.Range("E17").Value = 10
$ws.Range("E18").Value = 10
$ws.Range("E19").Value = 10
$ws.Range("E20").Value = 10
$ws.Range("E22").Value = 60

# --- ELE_NEW_WIND-OFF upper bound total for 2050 (row 31) ---
$ws.Range("E31").Value = 13.9

# --- ELE_NEW_PV_GRND upper-bound block (rows 33-38, 40) ---
$ws.Range("E33").Value = 10
$ws.Range("E34").Value = 10
$ws.Range("E35").Value = 10
$ws.Range("E36").Value = 10
$ws.Range("E37").Value = 10
$ws.Range("E38").Value = 10
$ws.Range("E40").Value = 60

# --- Restore sheet view state (zoom + active selection) ---
$ws.Range("E41").Select()
$excel.ActiveWindow.Zoom = 97
$excel.ActiveWindow.ScrollRow = 24
$excel.ActiveWindow.ScrollColumn = 1
